$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (46075 -> 46076) for every data row (rows 2 through 547).
$ws.Range("C2:C547").Value = 46076
